$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.021814
$ws.Range("H2").Value = 0.065442
$ws.Range("I2").Value = 0.1008129179549036
$ws.Range("J2").Value = 0.1008129179549036
$ws.Range("Q2").Value = 0.2512368843053333
$ws.Range("R2").Value = 2.261131958748
$ws.Range("S2").Value = 0.1008129179549036
$ws.Range("T2").Value = 0.1008129179549036

# Row 3 updates
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.194567
$ws.Range("H3").Value = 0.583701
$ws.Range("I3").Value = 0.8991870820450963
$ws.Range("J3").Value = 0.8991870820450963
$ws.Range("Q3").Value = 2.240873148832667
$ws.Range("R3").Value = 20.167858339494
$ws.Range("S3").Value = 0.8991870820450963
$ws.Range("T3").Value = 0.8991870820450963
